# Enable suppression of age calculation plots.
#
# On the "constants" sheet, insert a new boolean output-toggle row just above
# the existing "output_comorbidity_fractions" row (old row 93), labelled
# "output_age_calculations" and defaulting to FALSE. Everything at/below the
# old row 93 shifts down by one row (old row 117 -> new row 118).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Copy row 93 (an existing "output_*" boolean-toggle row that already has the
# exact style/number-format/data-validation this new row needs) and insert
# the copy above itself. This shifts row 93 and everything below it down by
# one row, while the newly-inserted row 93 inherits matching formatting.
$ws.Rows("93:93").Copy() | Out-Null
$ws.Rows("93:93").Insert(-4121) | Out-Null   # -4121 = xlShiftDown
$excel.CutCopyMode = 0

# Retarget the new row's label/value: a new output toggle, default off.
$ws.Range("A93").Value = "output_age_calculations"
$ws.Range("B93").Value = $false

# Match the author's recorded view state after the edit.
$ws.Activate()
$ws.Range("A94").Select() | Out-Null
